$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(13, 1).Value = "com.singleton.strechy"
$ws.Cells.Item(13, 2).Value = "stretchy"
$ws.Cells.Item(13, 3).Value = "ronoren61@gmail.com"
$ws.Cells.Item(13, 4).Value = "nitanoren23@gmail.com"
$ws.Cells.Item(13, 5).Value = "27/5/2019 15:59"
$ws.Cells.Item(13, 6).Value = "I can play this game for days months and even years. Never enough. can play offline cars too. Great taxi and cars game."

$ws.Hyperlinks.Add($ws.Cells.Item(13, 3), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13, 4), "mailto:nitanoren23@gmail.com", "", "", "nitanoren23@gmail.com") | Out-Null

$ws.Cells.Item(13, 1).Font.Name = "Mangal"
$ws.Cells.Item(13, 1).Font.Size = 10

$ws.Cells.Item(13, 6).Font.Name = "Mangal"
$ws.Cells.Item(13, 6).Font.Size = 10

$ws.Cells.Item(13, 3).Font.Name = "Calibri"
$ws.Cells.Item(13, 3).Font.Size = 11
$ws.Cells.Item(13, 3).Font.Underline = $false
$ws.Cells.Item(13, 3).Font.Color = 0
$ws.Cells.Item(13, 3).HorizontalAlignment = -4108

$ws.Cells.Item(13, 4).Font.Name = "Calibri"
$ws.Cells.Item(13, 4).Font.Size = 11
$ws.Cells.Item(13, 4).Font.Underline = $false
$ws.Cells.Item(13, 4).Font.Color = 0
$ws.Cells.Item(13, 4).HorizontalAlignment = -4108

$ws.Range("F13").Select()
